# "Loan RBI, Variable Instalments"
# The Repayment schedule sheet gets a new blank column inserted before column
# N ("Late"/"heading"/"Outstanding" shift right by one to O/P/Q), and becomes
# the active/selected worksheet tab (selection on J20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N; existing N/O/P shift to O/P/Q.
$ws.Columns("N:N").Insert() | Out-Null

# The new column N keeps the same on-screen width as column M ("In Advance"),
# i.e. a stored column width of 11 characters.
$ws.Columns("N:N").ColumnWidth = 11 - 5/6

# Make "Repayment schedule" the active sheet/tab with the same selection.
$ws.Activate() | Out-Null
$ws.Range("J20").Select() | Out-Null
